# Add the "26 deg samples" rows (17-26) to the CRM accuracy data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: a full data row (date, CRM value, batch value, % off, batch #, flags) ---
# Copy A16's style (date format) down to A17 before writing its value.
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value2 = 43194

$ws.Range("B17").Value2 = 2222.7716464551399

# --- Rows 17-26: batch value / % off / batch # populated for every new sample.
# Rows 18-26 have no CRM / batch-value reading (only C, D, E get filled in).
$ws.Range("C17:C26").Value2 = 2207.0300000000002
$ws.Range("D17:D26").Formula = "=100*(B17-C17)/C17"
$ws.Range("E17:E26").Value2 = 169

$ws.Range("F17").Value = "With Junk"
$ws.Range("G17").Value = "end of sample"

# Match the author's final selection/cursor position recorded in the diff.
$ws.Range("A18").Select()
